$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Database string only replacement: update the phone number and cashapp
# values for the Erik Bridges row (row 2) with their new shared-string text.
$ws.Range("B2").Value = " 555-555-555"
$ws.Range("C2").Value = " `$app"
